$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Session" to "Anatomy"
$ws.Name = "Anatomy"

# Update row 4 Student ID value - keep it stored as text (matching the rest
# of the "Student ID" column, which is all numeric-looking text) instead of
# letting Excel auto-convert the digit string into a real number.
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "235000"
$ws.Cells.Item(4, 1).Style = "Normal"

# Delete rows 5 and 6 (the trailing two data rows) - delete from bottom up
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
